# TS5-Preskladnenie tovaru.xlsx
# The "Testovane data" column (column D, between "Postup kroku:" and
# "Dodatocne informacie:") is removed from the test-scenario table, and
# the sheet's print orientation is switched from portrait to landscape.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the whole "Testované dáta" column - its header lived in D5 and the
# rest of the column was blank; deleting it shifts the following columns
# (Dodatočné informácie / Očakávaný výsledok) one slot to the left.
$ws.Columns("D").Delete()

# Keep the picture logo anchored at the same visual spot (top of what is
# now column E) instead of letting it trail off with the shifted columns.
$shp = $ws.Shapes.Item(1)
$shp.Left = 618.4658203125
$shp.Top = 4.2

# Switch the page to landscape printing.
$ws.PageSetup.Orientation = 2

# Park the selection back on a sensible cell.
$ws.Range("H6").Select() | Out-Null
